# Updated cryptos list on Wed Aug  9 19:02:59 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.505.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -1.14%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.852.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -0.34%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''0.9992'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.08%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''243.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -1.26%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.6529'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +2.55%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -0.06%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''47.96'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +2.67%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.07503'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.28%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.2982'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -0.58%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''24.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.90%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.07633'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -0.57%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''1.850.84'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -1.51%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''5.019'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.73%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''0.6850'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -0.57%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''83.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -0.81%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''0.000009525'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.98%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''6.119'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = '''29.538.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -0.88%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''2.111.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.33%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''236.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -1.41%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '''  -0.53%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  -0.11%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''7.691'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +4.66%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''1.001'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.04%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''157.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -1.21%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''0.1420'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -0.15%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''8.500'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -0.65%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  -0.96%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''0.06038'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -0.03%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''1.484'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -1.68%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''1.242'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -2.01%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''4.141'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -0.07%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''4.074'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -1.48%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''1.180'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +1.74%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''1.854'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -0.53%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.7243'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -0.97%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''2.599'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -0.77%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''2.805'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -2.21%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''0.01784'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -0.32%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''1.202.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -1.73%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''6.252'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -1.26%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.9080'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -1.60%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''0.9995'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = '''2.018.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +0.09%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''102.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -0.35%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''66.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.07%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = '''BabyDogeCoin'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''0.00000000124'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +0.72%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = '''Aptos'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = '''7.404'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +10.39%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.4058'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -0.80%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''9.076'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -2.47%  '
$ws.Range('E51').Style = 'Normal'
